# Applies the cryptos-list refresh described by the commit diff.
# Each entry: cell reference, new value, and whether the value must be
# forced to Text (because it would otherwise be auto-parsed as a number,
# e.g. "1.00" -> 1 or "0.0000118" -> 1.18E-05).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '61.469.05'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  +7.89%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '3.401.20'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  +4.73%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  +0.02%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '411.89'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  +3.84%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '121.79'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  +12.64%  '; ForceText = $false }
    @{ Cell = 'D7'; Value = '3.397.98'; ForceText = $false }
    @{ Cell = 'E7'; Value = '  +4.74%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.577'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  -0.84%  '; ForceText = $false }
    @{ Cell = 'E9'; Value = '  +0.05%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.641'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  +3.51%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.112'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  +16.67%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '41.42'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  +5.23%  '; ForceText = $false }
    @{ Cell = 'E13'; Value = '  -0.74%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '3.943.71'; ForceText = $false }
    @{ Cell = 'E14'; Value = '  +5.04%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '8.39'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  +0.90%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '19.51'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  +3.16%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '3.434.71'; ForceText = $false }
    @{ Cell = 'E17'; Value = '  +6.23%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '61.395.31'; ForceText = $false }
    @{ Cell = 'E18'; Value = '  +8.03%  '; ForceText = $false }
    @{ Cell = 'E19'; Value = '  -0.75%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '10.84'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -2.25%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '0.0000118'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  +5.95%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '3.33'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -0.20%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '12.84'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  -1.40%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '298.95'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  +1.81%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '76.05'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  +2.10%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '3.12'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -1.59%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '30.77'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  +9.40%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '8.18'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  +13.22%  '; ForceText = $false }
    @{ Cell = 'B29'; Value = 'LEO'; ForceText = $false }
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; ForceText = $false }
    @{ Cell = 'D29'; Value = '4.27'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  -2.04%  '; ForceText = $false }
    @{ Cell = 'B30'; Value = 'Filecoin'; ForceText = $false }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; ForceText = $false }
    @{ Cell = 'D30'; Value = '7.65'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  -2.56%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '0.170'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  +0.28%  '; ForceText = $false }
    @{ Cell = 'E32'; Value = '  +4.90%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '42.34'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  +2.04%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '11.40'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  +2.14%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '1.00'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  +0.13%  '; ForceText = $false }
    @{ Cell = 'E36'; Value = '  +17.43%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '0.0480'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  -1.06%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '52.26'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +2.02%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '3.53'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  +2.11%  '; ForceText = $false }
    @{ Cell = 'E40'; Value = '  -0.14%  '; ForceText = $false }
    @{ Cell = 'E41'; Value = '  +1.58%  '; ForceText = $false }
    @{ Cell = 'E43'; Value = '  +0.33%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '133.97'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  -2.00%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '17.30'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  +2.83%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '3.93'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  +2.30%  '; ForceText = $false }
    @{ Cell = 'B48'; Value = 'WEMIXToken'; ForceText = $false }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; ForceText = $false }
    @{ Cell = 'D48'; Value = '2.20'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  -3.14%  '; ForceText = $false }
    @{ Cell = 'B49'; Value = 'EnergySwap'; ForceText = $false }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText = $false }
    @{ Cell = 'D49'; Value = '21.92'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  -2.27%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '2.200.05'; ForceText = $false }
    @{ Cell = 'E50'; Value = '  +2.15%  '; ForceText = $false }
    @{ Cell = 'D51'; Value = '3.742.75'; ForceText = $false }
    @{ Cell = 'E51'; Value = '  +5.04%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $value = $u.Value
    if ($u.ForceText) {
        # Prefix with an apostrophe so Excel stores the literal digits/dots
        # as text instead of silently coercing them to a Double.
        $value = "'" + $value
    }
    $ws.Range($u.Cell).Value = $value
}
